$wb = $excel.ActiveWorkbook

# --- "Top Gainers" sheet: rows 40-44 shift up by one, row 44 gets new INDOTHAI data ---
$gainers = $wb.Worksheets.Item("Top Gainers")

$gainers.Cells.Item(40, 2).Value = "SURYAROSNI"
$gainers.Cells.Item(40, 3).Value = 4.9831
$gainers.Cells.Item(40, 4).Value = 11.386
$gainers.Cells.Item(40, 5).Value = 3.0213

$gainers.Cells.Item(41, 2).Value = "PDSL"
$gainers.Cells.Item(41, 3).Value = 4.9424
$gainers.Cells.Item(41, 4).Value = 2.96
$gainers.Cells.Item(41, 5).Value = 8.7852

$gainers.Cells.Item(42, 2).Value = "BIL"
$gainers.Cells.Item(42, 3).Value = 4.9242
$gainers.Cells.Item(42, 4).Value = 9.7065
$gainers.Cells.Item(42, 5).Value = 0.2135

$gainers.Cells.Item(43, 2).Value = "HITECHGEAR"
$gainers.Cells.Item(43, 3).Value = 4.8651
$gainers.Cells.Item(43, 4).Value = 2.1287
$gainers.Cells.Item(43, 5).Value = 10.9905

$gainers.Cells.Item(44, 2).Value = "INDOTHAI"
$gainers.Cells.Item(44, 3).Value = 4.8064
$gainers.Cells.Item(44, 4).Value = 4.5349
$gainers.Cells.Item(44, 5).Value = 43.748

# --- "Top Losers" sheet: rows 47-49 shift up by one, row 49 gets new DIGITIDE data ---
$losers = $wb.Worksheets.Item("Top Losers")

$losers.Cells.Item(47, 2).Value = "POLICYBZR"
$losers.Cells.Item(47, 3).Value = -2.907
$losers.Cells.Item(47, 4).Value = 2.2365
$losers.Cells.Item(47, 5).Value = 1.2573

$losers.Cells.Item(48, 2).Value = "BOSCHLTD"
$losers.Cells.Item(48, 3).Value = -2.9061
$losers.Cells.Item(48, 4).Value = -3.0193
$losers.Cells.Item(48, 5).Value = -1.9006

$losers.Cells.Item(49, 2).Value = "DIGITIDE"
$losers.Cells.Item(49, 3).Value = -2.8795
$losers.Cells.Item(49, 4).Value = 3.2317
$losers.Cells.Item(49, 5).Value = 6.2968

# --- "1 Month Performance" sheet: INDOTHAI's monthly change value updated ---
$perf = $wb.Worksheets.Item("1 Month Performance")
$perf.Cells.Item(2, 3).Value = 106.9036
